$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 568.7049
$ws.Range("J17").Value = 581.2373
$ws.Range("L17").Value = 1743.7119
$ws.Range("N17").Value = -2079.7119
$ws.Range("H98").Value = 667.8095
$ws.Range("I98").Value = 516.7222
$ws.Range("K98").Value = 516.7222
$ws.Range("M98").Value = 981.2778
$ws.Range("H107").Value = 144.08333
$ws.Range("I107").Value = 160.66667
$ws.Range("J107").Value = 94.333336
$ws.Range("K107").Value = 160.66667
$ws.Range("L107").Value = 94.333336
$ws.Range("M107").Value = 1759.33333
$ws.Range("N107").Value = -3934.333336
$ws.Range("H122").Value = 667.8095
$ws.Range("I122").Value = 516.7222
$ws.Range("K122").Value = 1550.1666
$ws.Range("M122").Value = 899.8334
$ws.Range("H141").Value = 7310.1816
$ws.Range("I141").Value = 6301.625
$ws.Range("K141").Value = 18904.875
$ws.Range("M141").Value = -13724.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 4504
$ws.Range("I22").Value = 5990
$ws.Range("J22").Value = 3018
$ws.Range("K22").Value = 5990
$ws.Range("L22").Value = 3018
$ws.Range("M22").Value = -5691
$ws.Range("N22").Value = -3616
$ws.Range("H32").Value = 5766.875
$ws.Range("I32").Value = 2358.7234
$ws.Range("K32").Value = 2358.7234
$ws.Range("M32").Value = -2071.7234
$ws.Range("H97").Value = 704.7692
$ws.Range("I97").Value = 637.5454999999999
$ws.Range("J97").Value = 1074.5
$ws.Range("K97").Value = 637.5454999999999
$ws.Range("L97").Value = 1074.5
$ws.Range("M97").Value = -141.5454999999999
$ws.Range("N97").Value = -2066.5
$ws.Range("H110").Value = 1083.6364
$ws.Range("I110").Value = 1092
$ws.Range("K110").Value = 1092
$ws.Range("M110").Value = 953
$ws.Range("H122").Value = 2359.7334
$ws.Range("I122").Value = 2322.8462
$ws.Range("K122").Value = 6968.5386
$ws.Range("M122").Value = -4518.5386
$ws.Range("H132").Value = 1440.619
$ws.Range("I132").Value = 1322.0555
$ws.Range("K132").Value = 3966.1665
$ws.Range("M132").Value = -1436.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2518.1428
$ws.Range("I16").Value = 2415.4
$ws.Range("K16").Value = 2415.4
$ws.Range("M16").Value = -2128.4
$ws.Range("H31").Value = 3206.5715
$ws.Range("I31").Value = 2061.5
$ws.Range("K31").Value = 2061.5
$ws.Range("M31").Value = -1766.5
$ws.Range("H34").Value = 3206.5715
$ws.Range("I34").Value = 2061.5
$ws.Range("K34").Value = 2061.5
$ws.Range("M34").Value = -1859.5
$ws.Range("H59").Value = 92498
$ws.Range("J59").Value = 169996
$ws.Range("L59").Value = 169996
$ws.Range("N59").Value = -172286
$ws.Range("H105").Value = 4999.8
$ws.Range("I105").Value = 1666.3334
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 1666.3334
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = 80.66660000000002
$ws.Range("N105").Value = -13494
$ws.Range("H113").Value = 2518.1428
$ws.Range("I113").Value = 2415.4
$ws.Range("K113").Value = 2415.4
$ws.Range("M113").Value = -245.4000000000001
$ws.Range("H122").Value = 1832.0741
$ws.Range("I122").Value = 1540.3684
$ws.Range("K122").Value = 4621.1052
$ws.Range("M122").Value = -2171.1052
$ws.Range("H132").Value = 1676.4546
$ws.Range("I132").Value = 1519.1
$ws.Range("K132").Value = 4557.299999999999
$ws.Range("M132").Value = -2027.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 167833
$ws.Range("J113").Value = 334199.34
$ws.Range("L113").Value = 1002598.02
$ws.Range("N113").Value = -1006938.02
$ws.Range("H136").Value = 1466.5625
$ws.Range("I136").Value = 1466.5625
$ws.Range("K136").Value = 4399.6875
$ws.Range("M136").Value = 700.3125
$ws.Range("H139").Value = 3469.5557
$ws.Range("I139").Value = 2413.4666
$ws.Range("J139").Value = 8750
$ws.Range("K139").Value = 7240.399800000001
$ws.Range("L139").Value = 26250
$ws.Range("M139").Value = -2100.399800000001
$ws.Range("N139").Value = -36530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 18394
$ws.Range("H80").Value = 100003120
$ws.Range("I80").Value = 200002030
$ws.Range("J80").Value = 4199.2
$ws.Range("K80").Value = 200002030
$ws.Range("L80").Value = 4199.2
$ws.Range("M80").Value = -200001032
$ws.Range("N80").Value = -6195.2
$ws.Range("H83").Value = 100003120
$ws.Range("I83").Value = 200002030
$ws.Range("J83").Value = 4199.2
$ws.Range("K83").Value = 1000010150
$ws.Range("L83").Value = 20996
$ws.Range("M83").Value = -1000005158
$ws.Range("N83").Value = -30980
$ws.Range("H93").Value = 18695.5
$ws.Range("J93").Value = 18695.5
$ws.Range("L93").Value = 18695.5
$ws.Range("N93").Value = -22439.5
$ws.Range("H102").Value = 1816.7222
$ws.Range("I102").Value = 1753.4667
$ws.Range("J102").Value = 2133
$ws.Range("K102").Value = 1753.4667
$ws.Range("L102").Value = 2133
$ws.Range("M102").Value = -131.4666999999999
$ws.Range("N102").Value = -5377
$ws.Range("H122").Value = 1120268.2
$ws.Range("I122").Value = 1437343.8
$ws.Range("J122").Value = 10504
$ws.Range("K122").Value = 4312031.4
$ws.Range("L122").Value = 31512
$ws.Range("M122").Value = -4309581.4
$ws.Range("N122").Value = -36412
$ws.Range("H123").Value = 51999
$ws.Range("J123").Value = 51999
$ws.Range("L123").Value = 51999
$ws.Range("N123").Value = -56899
$ws.Range("H126").Value = 3966.6956
$ws.Range("I126").Value = 2616.5
$ws.Range("K126").Value = 7849.5
$ws.Range("M126").Value = -5379.5
$ws.Range("H132").Value = 3818
$ws.Range("I132").Value = 3068.647
$ws.Range("J132").Value = 4667.2666
$ws.Range("K132").Value = 9205.940999999999
$ws.Range("L132").Value = 14001.7998
$ws.Range("M132").Value = -6675.940999999999
$ws.Range("N132").Value = -19061.7998
$ws.Range("H139").Value = 60296
$ws.Range("I139").Value = 60296
$ws.Range("K139").Value = 60296
$ws.Range("M139").Value = -55156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -888
$ws.Range("H7").Value = 5303.125
$ws.Range("I7").Value = 3700
$ws.Range("K7").Value = 3700
$ws.Range("M7").Value = -3588
$ws.Range("H22").Value = 1394
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 1394
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H40").Value = 6483858
$ws.Range("I40").Value = 2592.5
$ws.Range("J40").Value = 77777780
$ws.Range("K40").Value = 2592.5
$ws.Range("L40").Value = 77777780
$ws.Range("M40").Value = -2456.5
$ws.Range("N40").Value = -77778052
$ws.Range("H110").Value = 19527.8
$ws.Range("J110").Value = 19527.8
$ws.Range("L110").Value = 19527.8
$ws.Range("N110").Value = -27707.8
$ws.Range("H122").Value = 50003644
$ws.Range("I122").Value = 76926420
$ws.Range("K122").Value = 230779260
$ws.Range("M122").Value = -230776810
$ws.Range("H126").Value = 5303.125
$ws.Range("I126").Value = 3700
$ws.Range("K126").Value = 11100
$ws.Range("M126").Value = -8630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 28974.354
$ws.Range("J4").Value = 17813
$ws.Range("L4").Value = 17813
$ws.Range("N4").Value = -18039
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H27").Value = 50619.5
$ws.Range("J27").Value = 50619.5
$ws.Range("L27").Value = 50619.5
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = -50757.5
$ws.Range("H54").Value = 31250
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 31250
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 31250
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -32290
$ws.Range("H70").Value = 53000
$ws.Range("J70").Value = 53000
$ws.Range("L70").Value = 53000
$ws.Range("N70").Value = -53630
$ws.Range("H73").Value = 53000
$ws.Range("J73").Value = 53000
$ws.Range("L73").Value = 53000
$ws.Range("N73").Value = -55184
$ws.Range("H122").Value = 2143.9285
$ws.Range("I122").Value = 1721.5
$ws.Range("K122").Value = 5164.5
$ws.Range("M122").Value = -2714.5
$ws.Range("H126").Value = 12749.5
$ws.Range("I126").Value = 15713.429
$ws.Range("K126").Value = 47140.287
$ws.Range("M126").Value = -44670.287
$ws.Range("H132").Value = 1088.8966
$ws.Range("I132").Value = 928.8148
$ws.Range("K132").Value = 2786.4444
$ws.Range("M132").Value = -256.4443999999999
